$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 26 (last row of old data), since the new layout
# only needs rows 1-6.
$ws.Range("A7:A26").EntireRow.Delete() | Out-Null

# Update A2:A6 with the condensed per-card tuple strings.
$ws.Range("A2").Value = "('Cathedral of War', ['Land', 'Cathedral of War enters the battlefield tapped.', 'Exalted (Whenever a creature you control attacks alone, that creature gets +1/+1 until end of turn.)', '{T}: Add {C}.'])"

$ws.Range("A3").Value = "('Magmaquake', ['{X}{R}{R}', 'Instant', 'Magmaquake deals X damage to each creature without flying and each planeswalker.'])"

$ws.Range("A4").Value = "('Mwonvuli Beast Tracker', ['{1}{G}{G}', 'Creature — Human Scout', 'When Mwonvuli Beast Tracker enters the battlefield, search your library for a creature card with deathtouch, hexproof, reach, or trample and reveal it. Shuffle your library and put that card on top of it.', '2/1'])"

$ws.Range("A5").Value = "('Staff of Nin', ['{6}', 'Artifact', 'At the beginning of your upkeep, draw a card.', '{T}: Staff of Nin deals 1 damage to any target.'])"

$ws.Range("A6").Value = "('Xathrid Gorgon', ['{5}{B}', 'Creature — Gorgon', 'Deathtouch (Any amount of damage this deals to a creature is enough to destroy it.)', '{2}{B}, {T}: Put a petrification counter on target creature. It gains defender and becomes a colorless artifact in addition to its other types. Its activated abilities can’t be activated. (A creature with defender can’t attack.)', '3/6'])"
